$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "does"
$ws.Range("D3").Value = "this"
$ws.Range("H3").Value = 19.0
$ws.Range("G8").Value = "location 19"
$ws.Range("H8").Value = "name 19"
